# Update the date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-02-01 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-02-02 Friday", 2)

# Update the division problems in the single 20-row x 5-column table.
# Only rows 1, 5, 9, 13, 17 (1-based) carry real content; the rest are
# blank spacer rows. Values are addressed positionally (row, column)
# rather than by text search, because several values repeat across
# different cells before/after the edit (e.g. "423÷4=").
$t = $word.ActiveDocument.Tables.Item(1)

$values = @{
    1  = @("696÷2=", "154÷4=", "423÷4=", "212÷5=", "835÷4=")
    5  = @("699÷5=", "293÷5=", "497÷2=", "269÷7=", "354÷5=")
    9  = @("331÷8=", "539÷3=", "187÷3=", "260÷5=", "780÷3=")
    13 = @("321÷2=", "605÷4=", "114÷9=", "586÷6=", "498÷2=")
    17 = @("240÷3=", "296÷9=", "418÷3=", "978÷7=", "646÷7=")
}

foreach ($rowIndex in $values.Keys) {
    $rowValues = $values[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $rowValues[$col - 1]
    }
}
